$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current_investments")

$ws.Range("A7").Value = 12
$ws.Range("A8").Value = 14
